$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4702473333333333
$ws.Range("N2").Value = 1.410742
$ws.Range("O2").Value = 0.009034922268422819
$ws.Range("P2").Value = 0.009034922268422819
$ws.Range("Q2").Value = 9.962659847250888
$ws.Range("R2").Value = 89.663938625258
$ws.Range("S2").Value = 0.0007919237303056096
$ws.Range("T2").Value = 0.0007919237303056096
$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("N3").Value = 0.9584440000000001
$ws.Range("O3").Value = 0.006138235792679485
$ws.Range("P3").Value = 0.006138235792679485
$ws.Range("Q3").Value = 6.768531421506222
$ws.Range("R3").Value = 60.91678279355601
$ws.Range("S3").Value = 0.0005380250589895458
$ws.Range("T3").Value = 0.0005380250589895458
$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 1.047307
$ws.Range("N4").Value = 3.141921
$ws.Range("O4").Value = 0.02012204358311108
$ws.Range("P4").Value = 0.02012204358311108
$ws.Range("Q4").Value = 22.18824575289767
$ws.Range("R4").Value = 199.694211776079
$ws.Range("S4").Value = 0.00176372561293669
$ws.Range("T4").Value = 0.00176372561293669
$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("M5").Value = 50.21070966666667
$ws.Range("N5").Value = 150.632129
$ws.Range("O5").Value = 0.9647047983557866
$ws.Range("P5").Value = 0.9647047983557866
$ws.Range("Q5").Value = 1063.764078261097
$ws.Range("R5").Value = 9573.876704349872
$ws.Range("S5").Value = 0.08455774160091344
$ws.Range("T5").Value = 0.08455774160091344
$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4702473333333333
$ws.Range("N6").Value = 1.410742
$ws.Range("O6").Value = 0.009034922268422819
$ws.Range("P6").Value = 0.009034922268422819
$ws.Range("Q6").Value = 57.28594976403711
$ws.Range("R6").Value = 515.573547876334
$ws.Range("S6").Value = 0.004553613565734096
$ws.Range("T6").Value = 0.004553613565734096
$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("N7").Value = 0.9584440000000001
$ws.Range("O7").Value = 0.006138235792679485
$ws.Range("P7").Value = 0.006138235792679485
$ws.Range("S7").Value = 0.003093679496602817
$ws.Range("T7").Value = 0.003093679496602817
$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 1.047307
$ws.Range("N8").Value = 3.141921
$ws.Range("O8").Value = 0.02012204358311108
$ws.Range("P8").Value = 0.02012204358311108
$ws.Range("Q8").Value = 127.5838732869463
$ws.Range("R8").Value = 1148.254859582517
$ws.Range("S8").Value = 0.01014153834511543
$ws.Range("T8").Value = 0.01014153834511543
$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("M9").Value = 50.21070966666667
$ws.Range("N9").Value = 150.632129
$ws.Range("O9").Value = 0.9647047983557866
$ws.Range("P9").Value = 0.9647047983557866
$ws.Range("Q9").Value = 6116.713456283261
$ws.Range("R9").Value = 55050.42110654934
$ws.Range("S9").Value = 0.4862125789476804
$ws.Range("T9").Value = 0.4862125789476804
$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4702473333333333
$ws.Range("N10").Value = 1.410742
$ws.Range("O10").Value = 0.009034922268422819
$ws.Range("P10").Value = 0.009034922268422819
$ws.Range("Q10").Value = 17.49657827309711
$ws.Range("R10").Value = 157.469204457874
$ws.Range("S10").Value = 0.001390788779909873
$ws.Range("T10").Value = 0.001390788779909873
$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("N11").Value = 0.9584440000000001
$ws.Range("O11").Value = 0.006138235792679485
$ws.Range("P11").Value = 0.006138235792679485
$ws.Range("Q11").Value = 11.88700022142978
$ws.Range("R11").Value = 106.983001992868
$ws.Range("S11").Value = 0.0009448879819073496
$ws.Range("T11").Value = 0.0009448879819073499
$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 1.047307
$ws.Range("N12").Value = 3.141921
$ws.Range("O12").Value = 0.02012204358311108
$ws.Range("P12").Value = 0.02012204358311108
$ws.Range("Q12").Value = 38.96734250797633
$ws.Range("R12").Value = 350.706082571787
$ws.Range("S12").Value = 0.003097482370386086
$ws.Range("T12").Value = 0.003097482370386086
$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("M13").Value = 50.21070966666667
$ws.Range("N13").Value = 150.632129
$ws.Range("O13").Value = 0.9647047983557866
$ws.Range("P13").Value = 0.9647047983557866
$ws.Range("Q13").Value = 1868.199029653729
$ws.Range("R13").Value = 16813.79126688356
$ws.Range("S13").Value = 0.1485016217757298
$ws.Range("T13").Value = 0.1485016217757298
$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4702473333333333
$ws.Range("N14").Value = 1.410742
$ws.Range("O14").Value = 0.009034922268422819
$ws.Range("P14").Value = 0.009034922268422819
$ws.Range("Q14").Value = 28.91709278993266
$ws.Range("R14").Value = 260.2538351093939
$ws.Range("S14").Value = 0.00229859619247324
$ws.Range("T14").Value = 0.00229859619247324
$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("N15").Value = 0.9584440000000001
$ws.Range("O15").Value = 0.006138235792679485
$ws.Range("P15").Value = 0.006138235792679485
$ws.Range("Q15").Value = 19.64598351927867
$ws.Range("R15").Value = 176.813851673508
$ws.Range("S15").Value = 0.001561643255179772
$ws.Range("T15").Value = 0.001561643255179772
$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 1.047307
$ws.Range("N16").Value = 3.141921
$ws.Range("O16").Value = 0.02012204358311108
$ws.Range("P16").Value = 0.02012204358311108
$ws.Range("Q16").Value = 64.402435807283
$ws.Range("R16").Value = 579.6219222655469
$ws.Range("S16").Value = 0.005119297254672871
$ws.Range("T16").Value = 0.005119297254672871
$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("M17").Value = 50.21070966666667
$ws.Range("N17").Value = 150.632129
$ws.Range("O17").Value = 0.9647047983557866
$ws.Range("P17").Value = 0.9647047983557866
$ws.Range("Q17").Value = 3087.625697284201
$ws.Range("R17").Value = 27788.6312755578
$ws.Range("S17").Value = 0.2454328560314629
$ws.Range("T17").Value = 0.2454328560314629
